$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for B2:G11 (values shifted down by one row; new row of
# data inserted at the top; oldest row's data dropped from the bottom).
$data = @(
    @(-0.1510908837489806, 2.084195974016096, 20.51875844210171, 4.529763618788701, 4.628991593011878, 23),
    @(-0.2243426211202774, 2.382314769375817, 19.09061844583581, 4.369281227597487, 4.46620300197608, 22),
    @(-0.7120868115500602, 1.933618923903567, 10.88254014902561, 3.298869525917267, 3.300643293046703, 21),
    @(-0.3051420100610381, 1.578682252322016, 10.2251943824045, 3.197685785439917, 3.265784782946779, 20),
    @(-0.2804117665781, 1.683207935177004, 11.04924039551323, 3.324039770446983, 3.402952828034286, 19),
    @(-0.3173356620072107, 1.884213172660337, 10.94775440033029, 3.308739095234057, 3.388969517186058, 18),
    @(-0.2356007239738259, 1.940837229566887, 12.01945225291491, 3.466908169091721, 3.565345857874665, 17),
    @(-0.1538742921933148, 1.930586848953809, 10.98523007623265, 3.314397392623983, 3.41940389010477, 16),
    @(-0.1230455590706185, 1.907399787585887, 12.46868919605704, 3.531103113200893, 3.652819212576486, 15),
    @(-0.09178185018074302, 2.13619364861378, 13.88641155430669, 3.726447578365579, 3.865944277603689, 14)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
